# Apply updated vm_pu values for rows 2-25 (buses), columns B..F and I..N
# Column order corresponds to: B,C,D,E,F,I,J,K,L,M,N
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @(
    @(1.02,1.039937150655295,1.046812387593398,1.038392989847484,1.055606551737749,1.042818648575959,1.045026465710855,1.049576775548565,1.041181151534735,1.058346546582373,1.01887697790105),
    @(1.02,1.040940573279948,1.047621291418033,1.03924835236977,1.056648411584039,1.043127613897949,1.045675030215433,1.050197404736704,1.041846398266942,1.059201298722806,1.019096675242671),
    @(1.02,1.041590044616103,1.048144855207556,1.039802334546926,1.057323224069737,1.043326297892945,1.046094286002377,1.050598489875214,1.042276721408399,1.059754434157812,1.019238589531461),
    @(1.02,1.041863127463902,1.04836499624163,1.04003534879677,1.057607072439084,1.043409528105613,1.046270442726439,1.050766984965777,1.042457596078851,1.059986984504921,1.019298191477467),
    @(1.02,1.041908981910442,1.048401960907313,1.040074479934923,1.057654741027757,1.043423485425415,1.046300014429214,1.050795268914738,1.042487963768185,1.060026031429157,1.0193081954494),
    @(1.02,1.04159369338569,1.048147796607068,1.039805447623733,1.057327016248702,1.043327411184741,1.046096640204494,1.050600741791426,1.042279138395747,1.059757541459741,1.019239386167041),
    @(1.02,1.040276222524034,1.047085729089764,1.038681958448441,1.055958516185109,1.042923320822276,1.045245735446423,1.049786623984208,1.04140600256197,1.058635402582127,1.018951276127711),
    @(1.02,1.037956138961661,1.045215407713618,1.036706135653237,1.053552125203994,1.042201801364797,1.043743225041874,1.048348213819018,1.039866411645006,1.056658489072966,1.018441727972343),
    @(1.02,1.036410421379529,1.043969366301383,1.035391597784208,1.051951319141666,1.041714446985552,1.042739492365254,1.047386730841797,1.03883937151727,1.055340877168222,1.018100793036874),
    @(1.02,1.035741349823453,1.04343002535695,1.034823033335641,1.051258979268272,1.041501916816326,1.042304382531593,1.046969801367011,1.038394505593305,1.054770421953857,1.017952874329185),
    @(1.02,1.035492862201408,1.043229721361946,1.034611939951747,1.051001937088244,1.041422748160897,1.042142690613023,1.046814845354257,1.038229240440537,1.054558541904951,1.017897887147216),
    @(1.02,1.03554616203824,1.043272685823386,1.03465721580824,1.051057067903428,1.041439740306038,1.042177377364109,1.046848088003284,1.038264691370528,1.054603990325006,1.017909684053213),
    @(1.02,1.03572080903729,1.043413467524689,1.034805582312327,1.051237729555514,1.041495377309014,1.042291018508942,1.046956994484791,1.038380845172933,1.054752907626808,1.017948329959439),
    @(1.02,1.035828419505057,1.043500211939054,1.0348970086572,1.051349057565819,1.041529627230926,1.042361026940631,1.047024083471009,1.038452408445561,1.054844662151234,1.017972135208704),
    @(1.02,1.036454830404214,1.044005164926093,1.035429345063454,1.051997284773891,1.041728520277873,1.042768358946936,1.047414388454031,1.038868892670509,1.055378738116837,1.018110603789337),
    @(1.02,1.03684782424977,1.044321963504683,1.035763437265498,1.052404120404794,1.041852878541551,1.043023737466784,1.047659055871626,1.039130102216052,1.055713771467447,1.01819738355361),
    @(1.02,1.03707707351378,1.044506766314553,1.035958369210632,1.052641499800675,1.041925269625968,1.043172648424679,1.047801708218655,1.039282446743122,1.055909198384404,1.018247972568582),
    @(1.02,1.037155245418263,1.044569782592712,1.036024846394871,1.052722453433558,1.041949928515238,1.043223415226685,1.047850339127911,1.039334389801112,1.055975835192634,1.018265217333171),
    @(1.02,1.036805657390416,1.044287971987066,1.035727585955349,1.052360462610921,1.041839551071619,1.042996342628161,1.047632811371006,1.039102078408371,1.055677824766693,1.018188075816557),
    @(1.02,1.035669378872253,1.043372009954257,1.034761889397555,1.051184525783595,1.041478999820924,1.042257556029283,1.046924926716858,1.038346641380275,1.054709054846429,1.017936950904699),
    @(1.02,1.034955160030222,1.042796289403604,1.034155277880802,1.050445883435653,1.041251002461114,1.041792630848164,1.04647933196522,1.037871540290387,1.054100022830217,1.017778806631849),
    @(1.02,1.035333761515725,1.043101472352019,1.034476800769915,1.050837383694427,1.041371991741841,1.042039136226725,1.046715599272461,1.038123412314995,1.054422875230809,1.017862665721104),
    @(1.02,1.036824710709165,1.044303331230312,1.035743785428843,1.052380189444442,1.041845573629568,1.04300872132177,1.047644670310403,1.039114741205839,1.055694067509844,1.018192281669278),
    @(1.02,1.038555760367186,1.045698786162623,1.037216464921652,1.054173628878758,1.042389450580058,1.044132025139613,1.048720527912142,1.040264549302386,1.057169512969978,1.018573677392582)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $colLetters.Count; $j++) {
        $cellRef = "$($colLetters[$j])$rowNum"
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}

Write-Output "Updated $($data.Count) rows of vm_pu results"
